$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (GLD)
$ws.Range("D2").Value = 387.69
$ws.Range("E2").Value = 57
$ws.Range("F2").Value = 1.19
$ws.Range("N2").Value = 53.71147335634279

# Row 3 (NEM)
$ws.Range("D3").Value = 90.64
$ws.Range("E3").Value = 52.5
$ws.Range("F3").Value = 0.14
$ws.Range("H3").Value = 76
$ws.Range("I3").Value = 76
$ws.Range("K3").Value = 64.5
$ws.Range("N3").Value = 53.71147335634279

# Row 4 (Gold Feb 26)
$ws.Range("D4").Value = 4247
$ws.Range("E4").Value = 72
$ws.Range("F4").Value = 4.57
$ws.Range("H4").Value = 40
$ws.Range("I4").Value = 60
$ws.Range("K4").Value = 52.1
$ws.Range("N4").Value = 53.71147335634279
